# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310
#   *_new -> *_FV2404
# Then turn the used range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------------
$oldToFv2310 = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $oldToFv2310.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($oldToFv2310[$i])_FV2310"
}

# Column 11 is "diff" - unchanged.

for ($i = 0; $i -lt $oldToFv2310.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($oldToFv2310[$i])_FV2404"
}

# --- 2. Convert the used range into an Excel Table --------------------------
$tableRange = $ws.Range("A1:U91")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
